$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Tran Quoc Huy
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = 'Trần Quốc Huy'
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '03636'
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'Kế toán'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '15/03/2024'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = 7500000
$ws.Range("F5").Value = 2345678910
$ws.Range("G5").Value = 9000000
$ws.Range("H5").Value = 26
$ws.Range("I5").Value = 7500000
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 86538
$ws.Range("L5").Value = 400000
$ws.Range("M5").Value = 800000
$ws.Range("N5").Value = 300000
$ws.Range("O5").Value = 9086538
$ws.Range("P5").Value = 40000
$ws.Range("Q5").Value = 780000
$ws.Range("R5").Value = 820000
$ws.Range("S5").Value = 8266538
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 500000
$ws.Range("Y5").Value = 200000
$ws.Range("Z5").Value = 300000
$ws.Range("AA5").Value = 1
$ws.Range("AB5").Value = 400000
$ws.Range("AC5").Value = 0
$ws.Range("AD5").Value = 0
$ws.Range("AE5").Value = 26
$ws.Range("AF5").Value = 500000
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 0
$ws.Range("AK5").Value = 0
$ws.Range("AL5").Value = 4
$ws.Range("AM5").Value = 0
$ws.Range("AN5").Value = 8
$ws.Range("AO5").Value = 7500000
$ws.Range("AP5").Value = 0
$ws.Range("AQ5").Value = 300000
$ws.Range("AR5").Value = 200000
$ws.Range("AS5").Value = 0
$ws.Range("AT5").Value = 0
$ws.Range("AU5").Value = 200000
$ws.Range("AV5").Value = 300000
$ws.Range("AW5").Value = 750000
$ws.Range("AX5").Value = 350000
$ws.Range("AY5").Value = 0
$ws.Range("AZ5").Value = 0
$ws.Range("BA5").Value = 7500000
$ws.Range("BB5").Value = 0
$ws.Range("BC5").Value = 360000
$ws.Range("BD5").Value = 8
$ws.Range("BE5").NumberFormat = "@"
$ws.Range("BE5").Value = 'Khá'
$ws.Range("BE5").Style = "Normal"
$ws.Range("BF5").Value = -0.2
$ws.Range("BG5").Value = 120000

# Row 6 - Le Thi Mai
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = 'Lê Thị Mai'
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '01818'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'Nhân sự'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '01/06/2023'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = 6800000
$ws.Range("F6").Value = 9988776655
$ws.Range("G6").Value = 8500000
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = 6276923
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 350000
$ws.Range("M6").Value = 700000
$ws.Range("N6").Value = 300000
$ws.Range("O6").Value = 7626923
$ws.Range("P6").Value = 40000
$ws.Range("Q6").Value = 700000
$ws.Range("R6").Value = 740000
$ws.Range("S6").Value = 6886923
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 260000
$ws.Range("X6").Value = 300000
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 400000
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = 0
$ws.Range("AE6").Value = 24
$ws.Range("AF6").Value = 500000
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 1
$ws.Range("AK6").Value = 500000
$ws.Range("AL6").Value = 6
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = 6
$ws.Range("AO6").Value = 6800000
$ws.Range("AP6").Value = 0
$ws.Range("AQ6").Value = 300000
$ws.Range("AR6").Value = 0
$ws.Range("AS6").Value = 200000
$ws.Range("AT6").Value = 0
$ws.Range("AU6").Value = 200000
$ws.Range("AV6").Value = 0
$ws.Range("AW6").Value = 680000
$ws.Range("AX6").Value = 150000
$ws.Range("AY6").Value = 0
$ws.Range("AZ6").Value = 0
$ws.Range("BA6").Value = 6800000
$ws.Range("BB6").Value = 0
$ws.Range("BC6").Value = 3200000
$ws.Range("BD6").Value = 7
$ws.Range("BE6").NumberFormat = "@"
$ws.Range("BE6").Value = 'Khá'
$ws.Range("BE6").Style = "Normal"
$ws.Range("BF6").Value = -0.1
$ws.Range("BG6").Value = 180000

# Move the active selection to E12 (also clears the stale AZ1 scroll anchor)
$ws.Range("E12").Select() | Out-Null

